# emilia user story 3 return_inventory 19/12
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rebuild the inventory table: row 2 ("shirts", id 1) was removed, the
# remaining rows shifted up, the "shirt" row was corrected (size/amount/
# color), the "shirt"/"we"/"dc" row became "shorts", and a new "hats" row
# was appended.
$data = @(
    @(2, "jeans",  "s/m/l", 10, "black"),
    @(3, "shoes",  "s/m/l", 10, "black"),
    @(4, "coats",  "s/m/l", 10, "black"),
    @(6, "shirt",  "s/m/l", 6,  "white"),
    @(7, "shorts", "s/m/l", 10, "blue"),
    @(8, "hats",   "s/m/l", 10, "brown")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
